$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Density" header in column E, matching style of the other headers
$ws.Range("E1").Value = "Density"
$ws.Range("E1").Font.Bold = $true

# Density values for each material row
$density = 0.28356481
$ws.Range("E2").Value = $density
$ws.Range("E3").Value = $density
$ws.Range("E4").Value = $density
$ws.Range("E5").Value = $density
$ws.Range("E6").Value = $density
$ws.Range("E7").Value = $density
$ws.Range("E8").Value = $density

# Update selection to match the committed state
$ws.Range("E1:E8").Select()
